# Update cryptocurrency price (D) and 1h volume change (E) columns
# Values are forced as text (leading apostrophe) to avoid Excel
# auto-converting numeric-looking strings into numbers, then the
# style is reset to Normal so no extra formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.683.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.57%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.730.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.27%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'227.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.45%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.5435"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.98%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.2733"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.70%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.06677"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.58%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'21.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.34%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07767"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.91%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'4.683"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.18%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.737.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.26%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.969.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.37%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.5950"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +5.27%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0₅8380"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.30%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'68.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.77%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'27.690.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.67%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'224.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +17.16%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'4.804"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.21%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -0.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'10.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.36%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'6.199"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.16%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +0.02%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'147.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.87%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.734"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +13.47%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.1248"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.39%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'7.450"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.02%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'17.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.63%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.05658"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.05%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'1.311"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.41%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.661"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.54%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'3.496"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.02%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.671"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.58%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.9722"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.89%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'2.845"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.69%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.439"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.43%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.5972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.27%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01663"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.47%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'5.912"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.27%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.8588"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.80%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'1.047.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.71%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.07%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'101.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.29%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.874.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.32%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +11.15%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'59.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.39%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'8.250"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.40%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.4430"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.04%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.05326"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.58%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.57%  "
$ws.Range("E51").Style = "Normal"
